$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.682.36"

$ws.Range("D3").Value = "3.146.47"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.41%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.145.43"
$ws.Range("E8").Value = "  +0.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.25%  "

$ws.Range("E11").Value = "  -1.61%  "

$ws.Range("E12").Value = "  -1.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.71%  "

$ws.Range("D15").Value = "3.661.89"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").Value = "64.861.74"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("D17").Value = "3.149.38"
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.57%  "

$ws.Range("E19").Value = "  +0.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "500.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.710"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.37%  "

$ws.Range("E28").Value = "  -1.07%  "

$ws.Range("E29").Value = "  -0.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.80%  "

$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.30%  "

$ws.Range("E36").Value = "  -1.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0889"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "471.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0413"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("D42").Value = "3.000.33"
$ws.Range("E42").Value = "  -3.47%  "

$ws.Range("E43").Value = "  -3.62%  "

$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.70%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.279"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "27.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.67%  "

$ws.Range("D47").Value = "0.0₃0580"
$ws.Range("E47").Value = "  +1.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.113"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.70%  "

$ws.Range("E50").Value = "  -3.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.58%  "
